# Adding two new sections ("Facilities" and "equiptment") to the
# application tracker table, just above the "Applicant's Background..."
# row (the row that used to be row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the table body (old row 5),
# pushing all the existing section rows (and the blank rows further
# down the sheet) down by two rows.
$ws.Range("A5:D6").Insert(-4121)

# The structured table ("Table1") doesn't auto-grow when rows are
# inserted this way, so extend its range to include the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D17"))

# Fill in the two new rows. Write column A for both new rows before
# column B so new shared strings are appended in the same order as in
# the target workbook (Facilities, equiptment, ?).
$ws.Range("A5").Value2 = "Facilities"
$ws.Range("A6").Value2 = "equiptment"
$ws.Range("B5").Value2 = "?"
$ws.Range("B6").Value2 = "?"

# Match the widened "Section of Application" column (manual width of
# 67 characters, no more auto best-fit).
$ws.Columns("A").ColumnWidth = 66.16666666666667

# Update the selected cell shown when the workbook is reopened.
$ws.Range("B6").Select()
